$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '69.681.06'
$ws.Range('E2').Value = '  +0.36%  '

# Row 3
$ws.Range('D3').Value = '3.703.89'
$ws.Range('E3').Value = '  +0.63%  '

# Row 4
$ws.Range('E4').Value = '  -0.01%  '

# Row 5
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '675.20'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  -1.56%  '

# Row 6
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '161.91'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  +1.32%  '

# Row 7
$ws.Range('E7').Value = '  +0.08%  '

# Row 8
$ws.Range('E8').Value = '  +0.93%  '

# Row 9
$ws.Range('E9').Value = '  +0.99%  '

# Row 10
$ws.Range('E10').Value = '  +0.53%  '

# Row 11
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '0.445'
$c.Style = 'Normal'
$ws.Range('E11').Value = '  +2.50%  '

# Row 12
$ws.Range('E12').Value = '  +1.28%  '

# Row 13
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '32.90'
$c.Style = 'Normal'
$ws.Range('E13').Value = '  +2.01%  '

# Row 14
$ws.Range('D14').Value = '3.699.71'
$ws.Range('E14').Value = '  +0.29%  '

# Row 15
$ws.Range('D15').Value = '69.704.75'

# Row 17
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '16.18'
$c.Style = 'Normal'

# Row 18
$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '6.53'
$c.Style = 'Normal'
$ws.Range('E18').Value = '  +2.29%  '

# Row 19
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '474.45'
$c.Style = 'Normal'
$ws.Range('E19').Value = '  +0.83%  '

# Row 20
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '9.84'
$c.Style = 'Normal'
$ws.Range('E20').Value = '  -0.98%  '

# Row 21
$ws.Range('E21').Value = '  +0.82%  '

# Row 22
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '80.50'
$c.Style = 'Normal'
$ws.Range('E22').Value = '  +1.15%  '

# Row 23
$ws.Range('D23').Value = '3.850.40'
$ws.Range('E23').Value = '  +0.57%  '

# Row 24
$ws.Range('E24').Value = '  +2.92%  '

# Row 25
$ws.Range('E25').Value = '  -0.05%  '

# Row 26
$ws.Range('E26').Value = '  +0.43%  '

# Row 27
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '9.15'
$c.Style = 'Normal'
$ws.Range('E27').Value = '  -0.40%  '

# Row 28
$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '2.70'
$c.Style = 'Normal'
$ws.Range('E28').Value = '  -0.35%  '

# Row 29
$ws.Range('E29').Value = '  +1.21%  '

# Row 30
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '2.03'
$c.Style = 'Normal'
$ws.Range('E30').Value = '  +1.41%  '

# Row 31
$ws.Range('E31').Value = '  +0.81%  '

# Row 32
$ws.Range('E32').Value = '  +0.00%  '

# Row 33
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '26.91'
$c.Style = 'Normal'
$ws.Range('E33').Value = '  +0.26%  '

# Row 34
$ws.Range('E34').Value = '  +3.93%  '

# Row 35
$ws.Range('D35').Value = '3.692.56'
$ws.Range('E35').Value = '  +0.99%  '

# Row 36
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '8.52'
$c.Style = 'Normal'
$ws.Range('E36').Value = '  +4.06%  '

# Row 37
$ws.Range('E37').Value = '  +1.25%  '

# Row 38
$ws.Range('E38').Value = '  +0.00%  '

# Row 39
$ws.Range('E39').Value = '  -0.04%  '

# Row 40
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '2.22'
$c.Style = 'Normal'
$ws.Range('E40').Value = '  +0.46%  '

# Row 41
$ws.Range('E41').Value = '  +1.39%  '

# Row 42
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '170.85'
$c.Style = 'Normal'
$ws.Range('E42').Value = '  +3.36%  '

# Row 43
$ws.Range('E43').Value = '  +0.45%  '

# Row 44
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '46.99'
$c.Style = 'Normal'
$ws.Range('E44').Value = '  -1.07%  '

# Row 45
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '2.78'
$c.Style = 'Normal'
$ws.Range('E45').Value = '  +1.97%  '

# Row 46
$ws.Range('E46').Value = '  -1.21%  '

# Row 47
$ws.Range('B47').Value = 'SuiNetwork'
$ws.Range('C47').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '1.11'
$c.Style = 'Normal'
$ws.Range('E47').Value = '  -0.39%  '

# Row 48
$ws.Range('B48').Value = 'InjectiveProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '27.93'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  +1.06%  '

# Row 49
$ws.Range('E49').Value = '  -0.68%  '

# Row 50
$ws.Range('E50').Value = '  +1.99%  '

# Row 51
$ws.Range('E51').Value = '  +2.55%  '
